$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Thursday hours value for the week of row 7 (E7: 2.25 -> 4.5)
$ws.Range("E7").Value = 4.5

# Update the active selection to match the authored state (E8)
$ws.Range("E8").Select()

$wb.Save()
